$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking text values in columns B:E stay as text (matching original inlineStr type)
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "72.024.40"
$ws.Range("D3").Value = "3.622.83"
$ws.Range("E3").Value = "  +6.91%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "598.28"
$ws.Range("E5").Value = "  +1.75%  "
$ws.Range("D6").Value = "182.72"
$ws.Range("E6").Value = "  +1.55%  "
$ws.Range("D7").Value = "3.610.57"
$ws.Range("E7").Value = "  +6.77%  "
$ws.Range("D8").Value = "0.608"
$ws.Range("E8").Value = "  +2.02%  "
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("E10").Value = "  +6.15%  "
$ws.Range("D11").Value = "0.609"
$ws.Range("E11").Value = "  +3.08%  "
$ws.Range("D12").Value = "50.50"
$ws.Range("E12").Value = "  +4.03%  "
$ws.Range("D13").Value = "0.0000289"
$ws.Range("E13").Value = "  +2.72%  "
$ws.Range("D14").Value = "708.75"
$ws.Range("E14").Value = "  +4.13%  "
$ws.Range("D15").Value = "4.199.27"
$ws.Range("E15").Value = "  +6.81%  "
$ws.Range("D16").Value = "8.96"
$ws.Range("E16").Value = "  +3.86%  "
$ws.Range("D17").Value = "72.018.28"
$ws.Range("E17").Value = "  +3.77%  "
$ws.Range("D18").Value = "3.566.55"
$ws.Range("E18").Value = "  +5.27%  "
$ws.Range("D19").Value = "0.122"
$ws.Range("E19").Value = "  +1.69%  "
$ws.Range("D20").Value = "18.55"
$ws.Range("E20").Value = "  +4.85%  "
$ws.Range("E21").Value = "  +4.43%  "
$ws.Range("D22").Value = "0.934"
$ws.Range("E22").Value = "  +3.26%  "
$ws.Range("D23").Value = "5.84"
$ws.Range("E23").Value = "  +7.47%  "
$ws.Range("D24").Value = "17.86"
$ws.Range("E24").Value = "  +4.25%  "
$ws.Range("D25").Value = "105.72"
$ws.Range("E25").Value = "  +2.48%  "
$ws.Range("E26").Value = "  +2.87%  "
$ws.Range("D27").Value = "2.86"
$ws.Range("E27").Value = "  +4.80%  "
$ws.Range("E28").Value = "  +4.29%  "
$ws.Range("D29").Value = "35.77"
$ws.Range("E29").Value = "  +5.16%  "
$ws.Range("D30").Value = "9.12"
$ws.Range("E30").Value = "  +4.52%  "
$ws.Range("D31").Value = "7.44"
$ws.Range("E31").Value = "  +7.21%  "
$ws.Range("D32").Value = "4.20"
$ws.Range("E32").Value = "  +18.19%  "
$ws.Range("D33").Value = "592.17"
$ws.Range("E33").Value = "  +6.11%  "
$ws.Range("D34").Value = "11.40"
$ws.Range("E34").Value = "  +2.39%  "
$ws.Range("D36").Value = "59.84"
$ws.Range("E36").Value = "  +2.05%  "
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("D38").Value = "3.656.35"
$ws.Range("E38").Value = "  -0.38%  "
$ws.Range("E39").Value = "  +4.10%  "
$ws.Range("B40").Value = "InjectiveProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D40").Value = "36.09"
$ws.Range("E40").Value = "  +1.89%  "
$ws.Range("B41").Value = "PEPE"
$ws.Range("C41").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D41").Value = "0.0₃0775"
$ws.Range("E41").Value = "  +10.45%  "
$ws.Range("E42").Value = "  +5.10%  "
$ws.Range("E43").Value = "  +4.47%  "
$ws.Range("D44").Value = "0.0454"
$ws.Range("E44").Value = "  +7.37%  "
$ws.Range("D45").Value = "0.348"
$ws.Range("E45").Value = "  +2.12%  "
$ws.Range("D46").Value = "3.40"
$ws.Range("E46").Value = "  +2.43%  "
$ws.Range("D47").Value = "2.79"
$ws.Range("E47").Value = "  +4.18%  "
$ws.Range("E48").Value = "  +5.42%  "
$ws.Range("E49").Value = "  +2.10%  "
$ws.Range("D50").Value = "0.998"
$ws.Range("E50").Value = "  -0.34%  "
$ws.Range("D51").Value = "133.43"
$ws.Range("E51").Value = "  +0.17%  "
